$wb = $excel.ActiveWorkbook

# --- Rename Sheet7 to "byte vs pinyin" and update its selection ---
$sheet7 = $wb.Worksheets.Item("Sheet7")
$sheet7.Name = "byte vs pinyin"
$sheet7.Range("O27").Select()

# --- Add Sheet8 after the last sheet ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8 = $wb.Worksheets.Add([Type]::Missing, $last)
$ws8.Name = "Sheet8"

# Populate Sheet8 (order matters for shared-string allocation: "alex" must be entered first)
$ws8.Range("C11").Value = "alex"
$ws8.Range("C10").Value = "name"
$ws8.Range("D10").Value = "age"
$ws8.Range("E10").Value = "parent name"
$ws8.Range("F10").Value = "address"
$ws8.Range("G10").Value = "phone"
$ws8.Range("D11").Value = 1
$ws8.Range("F24").Value = "registration renewal?"
$ws8.Columns.Item(5).ColumnWidth = 17
$ws8.Columns.Item(6).ColumnWidth = 20.833333333333332
$ws8.Range("H10").Select()

# --- Add Sheet9 after Sheet8 ---
$last2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws9 = $wb.Worksheets.Add([Type]::Missing, $last2)
$ws9.Name = "Sheet9"

$ws9.Range("A1").Value = '2021-12-03 18:02:37 INFO  Benchmark_Timer - Begin run: msd collator and 1000000 elements with 10 runs'
$ws9.Range("A2").Value = 1193.3312632
$ws9.Range("A3").Value = '2021-12-03 18:02:52 INFO  Benchmark_Timer - Begin run: msd byte node and 1000000 elements with 10 runs'
$ws9.Range("A4").Value = 1002.4056785
$ws9.Range("A5").Value = '2021-12-03 18:03:04 INFO  Benchmark_Timer - Begin run: tim and 1000000 elements with 10 runs'
$ws9.Range("A6").Value = 1641.9662951999901
$ws9.Range("A7").Value = '2021-12-03 18:03:24 INFO  Benchmark_Timer - Begin run: quick and 1000000 elements with 10 runs'
$ws9.Range("A8").Value = 1622.5075199
$ws9.Range("A9").Value = '2021-12-03 18:03:43 INFO  Benchmark_Timer - Begin run: husky and 1000000 elements with 10 runs'
$ws9.Range("A10").Value = 1822.4092231
$ws9.Range("A11").Value = '2021-12-03 18:04:06 INFO  Benchmark_Timer - Begin run: lsd and 1000000 elements with 10 runs'
$ws9.Range("A12").Value = 2211.5141911000001
$ws9.Range("A13").Value = 2211.5141911000001
$ws9.Range("A14").Value = '2021-12-03 18:04:32 INFO  Benchmark_Timer - Begin run: lsd byte node and 1000000 elements with 10 runs'
$ws9.Range("A15").Value = 2907.9794781999999
$ws9.Range("A16").Value = '2021-12-03 18:05:07 INFO  Benchmark_Timer - Begin run: msd collator and 1500000 elements with 10 runs'
$ws9.Range("A17").Value = 4668.5997549000003
$ws9.Range("A18").Value = '2021-12-03 18:06:03 INFO  Benchmark_Timer - Begin run: msd byte node and 1500000 elements with 10 runs'
$ws9.Range("A19").Value = 3066.8674141000001
$ws9.Range("A20").Value = '2021-12-03 18:06:40 INFO  Benchmark_Timer - Begin run: tim and 1500000 elements with 10 runs'
$ws9.Range("A21").Value = 2464.0768904000001
$ws9.Range("A22").Value = '2021-12-03 18:07:10 INFO  Benchmark_Timer - Begin run: quick and 1500000 elements with 10 runs'
$ws9.Range("A23").Value = 2373.0892391000002
$ws9.Range("A24").Value = '2021-12-03 18:07:38 INFO  Benchmark_Timer - Begin run: husky and 1500000 elements with 10 runs'
$ws9.Range("A25").Value = 2705.8897711
$ws9.Range("A26").Value = '2021-12-03 18:08:10 INFO  Benchmark_Timer - Begin run: lsd and 1500000 elements with 10 runs'
$ws9.Range("A27").Value = 3514.3748728
$ws9.Range("A28").Value = 3514.3748728
$ws9.Range("A29").Value = '2021-12-03 18:08:53 INFO  Benchmark_Timer - Begin run: lsd byte node and 1500000 elements with 10 runs'
$ws9.Range("A30").Value = 4113.3844971999997
$ws9.Range("A31").Value = '2021-12-03 18:09:43 INFO  Benchmark_Timer - Begin run: msd collator and 2000000 elements with 10 runs'
$ws9.Range("A32").Value = 8087.30639229999
$ws9.Range("A33").Value = '2021-12-03 18:11:21 INFO  Benchmark_Timer - Begin run: msd byte node and 2000000 elements with 10 runs'
$ws9.Range("A34").Value = 5302.6826149999997
$ws9.Range("A35").Value = '2021-12-03 18:12:25 INFO  Benchmark_Timer - Begin run: tim and 2000000 elements with 10 runs'
$ws9.Range("A36").Value = 3537.6140154999998
$ws9.Range("A37").Value = '2021-12-03 18:13:08 INFO  Benchmark_Timer - Begin run: quick and 2000000 elements with 10 runs'
$ws9.Range("A38").Value = 3137.8088641999998
$ws9.Range("A39").Value = '2021-12-03 18:13:45 INFO  Benchmark_Timer - Begin run: husky and 2000000 elements with 10 runs'
$ws9.Range("A40").Value = 3403.7592008000001
$ws9.Range("A41").Value = '2021-12-03 18:14:26 INFO  Benchmark_Timer - Begin run: lsd and 2000000 elements with 10 runs'
$ws9.Range("A42").Value = 4338.0417101000003
$ws9.Range("A43").Value = 4338.0417101000003
$ws9.Range("A44").Value = '2021-12-03 18:15:17 INFO  Benchmark_Timer - Begin run: lsd byte node and 2000000 elements with 10 runs'
$ws9.Range("A45").Value = 5037.5307964000003
$ws9.Range("A46").Value = '2021-12-03 18:16:18 INFO  Benchmark_Timer - Begin run: msd collator and 3000000 elements with 10 runs'
$ws9.Range("A47").Value = 9824.8539834999992
$ws9.Range("A48").Value = '2021-12-03 18:18:18 INFO  Benchmark_Timer - Begin run: msd byte node and 3000000 elements with 10 runs'
$ws9.Range("A49").Value = 6293.7061626000004
$ws9.Range("A50").Value = '2021-12-03 18:19:34 INFO  Benchmark_Timer - Begin run: tim and 3000000 elements with 10 runs'
$ws9.Range("A51").Value = 5798.7852033999998
$ws9.Range("A52").Value = '2021-12-03 18:20:45 INFO  Benchmark_Timer - Begin run: quick and 3000000 elements with 10 runs'
$ws9.Range("A53").Value = 5743.1479866
$ws9.Range("A54").Value = '2021-12-03 18:21:54 INFO  Benchmark_Timer - Begin run: husky and 3000000 elements with 10 runs'
$ws9.Range("A55").Value = 5414.9443593999904
$ws9.Range("A56").Value = '2021-12-03 18:22:59 INFO  Benchmark_Timer - Begin run: lsd and 3000000 elements with 10 runs'
$ws9.Range("A57").Value = 6921.5136346999998
$ws9.Range("A58").Value = 6921.5136346999998
$ws9.Range("A59").Value = '2021-12-03 18:24:18 INFO  Benchmark_Timer - Begin run: lsd byte node and 3000000 elements with 10 runs'
$ws9.Range("A60").Value = 7047.3616068000001
$ws9.Range("A61").Value = '2021-12-03 18:25:51 INFO  Benchmark_Timer - Begin run: msd collator and 4000000 elements with 10 runs'
$ws9.Range("A62").Value = 10782.457198100001
$ws9.Range("A63").Value = '2021-12-03 18:28:03 INFO  Benchmark_Timer - Begin run: msd byte node and 4000000 elements with 10 runs'
$ws9.Range("A64").Value = 6810.1976133999997
$ws9.Range("A65").Value = '2021-12-03 18:29:24 INFO  Benchmark_Timer - Begin run: tim and 4000000 elements with 10 runs'
$ws9.Range("A66").Value = 7405.5059007999998
$ws9.Range("A67").Value = '2021-12-03 18:30:53 INFO  Benchmark_Timer - Begin run: quick and 4000000 elements with 10 runs'
$ws9.Range("A68").Value = 7640.6832867000003
$ws9.Range("A69").Value = '2021-12-03 18:32:25 INFO  Benchmark_Timer - Begin run: husky and 4000000 elements with 10 runs'
$ws9.Range("A70").Value = 7924.4500300999898
$ws9.Range("A71").Value = '2021-12-03 18:34:00 INFO  Benchmark_Timer - Begin run: lsd and 4000000 elements with 10 runs'
$ws9.Range("A72").Value = 9243.3009440999995
$ws9.Range("A73").Value = 9243.3009440999995
$ws9.Range("A74").Value = '2021-12-03 18:35:51 INFO  Benchmark_Timer - Begin run: lsd byte node and 4000000 elements with 10 runs'
$ws9.Range("A75").Value = 8465.1165122000002

$ws9.Range("I14").Select()
$ws9.Activate()
